# Add season-record columns (Wins / Losses / Ties) to the player table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new labeled columns AD/AE/AF ---
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the existing header formatting (bold, centered, bordered) by
# copying the format from the neighboring header cell (AC1) instead of
# re-building it from scratch, so the same style is reused.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# --- Data rows (2-48): same season record for every player on the roster ---
for ($r = 2; $r -le 48; $r++) {
    $ws.Cells.Item($r, 30).Value = 98   # AD -> Wins
    $ws.Cells.Item($r, 31).Value = 64   # AE -> Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF -> Ties
}
